$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix E40/F40: convert from text "21"/"20000" to numeric values,
# matching the filtered dataframe dtype used by the plotting code.
$ws.Range("E40").Value = 21
$ws.Range("F40").Value = 20000

# Append newly-logged algorithm runs (rows 41-61).
$ws.Range("A41").Value = "2025-06-25 17:52:11"
$ws.Range("B41").Value = "NaiveDFS"
$ws.Range("C41").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D41").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 20000
$ws.Range("G41").Value = 20.26472453499809
$ws.Range("H41").Value = 0.3377454089166349

$ws.Range("A42").Value = "2025-06-25 17:52:12"
$ws.Range("B42").Value = "Apriori"
$ws.Range("C42").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D42").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 20000
$ws.Range("G42").Value = 0.6112939700033166
$ws.Range("H42").Value = 0.01018823283338861

$ws.Range("A43").Value = "2025-06-25 17:52:18"
$ws.Range("B43").Value = "FP"
$ws.Range("C43").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D43").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 20000
$ws.Range("G43").Value = 5.423727711000538
$ws.Range("H43").Value = 0.09039546185000896

$ws.Range("A44").Value = "2025-06-25 20:08:03"
$ws.Range("B44").Value = "NaiveDFS"
$ws.Range("C44").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D44").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E44").Value = 25
$ws.Range("F44").Value = 20000
$ws.Range("G44").Value = 18.42766634200234
$ws.Range("H44").Value = 0.3071277723667057

$ws.Range("A45").Value = "2025-06-25 20:08:34"
$ws.Range("B45").Value = "NaiveDFS"
$ws.Range("C45").Value = "{'HoursComputer': '5 - 8 hours'}"
$ws.Range("D45").Value = "{'Gender': 'Male'}"
$ws.Range("E45").Value = 33
$ws.Range("F45").Value = 20000
$ws.Range("G45").Value = 30.30420459100424
$ws.Range("H45").Value = 0.5050700765167373

$ws.Range("A46").Value = "2025-06-25 20:09:14"
$ws.Range("B46").Value = "NaiveDFS"
$ws.Range("C46").Value = "{'HoursComputer': '1 - 4 hours', 'FormalEducation': 'Master’s degree (MA, MS, M.Eng., MBA, etc.)'}"
$ws.Range("D46").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E46").Value = 43
$ws.Range("F46").Value = 20000
$ws.Range("G46").Value = 40.33883208200132
$ws.Range("H46").Value = 0.6723138680333552

$ws.Range("A47").Value = "2025-06-25 20:23:09"
$ws.Range("B47").Value = "NaiveDFS"
$ws.Range("C47").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D47").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E47").Value = 21
$ws.Range("F47").Value = 20000
$ws.Range("G47").Value = 27.50110511200182
$ws.Range("H47").Value = 0.4583517518666971

$ws.Range("A48").Value = "2025-06-25 20:24:15"
$ws.Range("B48").Value = "NaiveDFS"
$ws.Range("C48").Value = "{'HoursComputer': '5 - 8 hours'}"
$ws.Range("D48").Value = "{'Gender': 'Male'}"
$ws.Range("E48").Value = 29
$ws.Range("F48").Value = 20000
$ws.Range("G48").Value = 64.83382585099753
$ws.Range("H48").Value = 1.080563764183292

$ws.Range("A49").Value = "2025-06-25 20:25:06"
$ws.Range("B49").Value = "NaiveDFS"
$ws.Range("C49").Value = "{'HoursComputer': '1 - 4 hours', 'FormalEducation': 'Master’s degree (MA, MS, M.Eng., MBA, etc.)'}"
$ws.Range("D49").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E49").Value = 21
$ws.Range("F49").Value = 20000
$ws.Range("G49").Value = 49.95099182800186
$ws.Range("H49").Value = 0.8325165304666977

$ws.Range("A50").Value = "2025-06-25 21:21:48"
$ws.Range("B50").Value = "Apriori"
$ws.Range("C50").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D50").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 20000
$ws.Range("G50").Value = 0.5756943699961994
$ws.Range("H50").Value = 0.009594906166603323

$ws.Range("A51").Value = "2025-06-25 21:21:48"
$ws.Range("B51").Value = "Apriori"
$ws.Range("C51").Value = "{'HoursComputer': '5 - 8 hours'}"
$ws.Range("D51").Value = "{'Gender': 'Male'}"
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 20000
$ws.Range("G51").Value = 0.6204752510020626
$ws.Range("H51").Value = 0.01034125418336771

$ws.Range("A52").Value = "2025-06-25 21:21:50"
$ws.Range("B52").Value = "Apriori"
$ws.Range("C52").Value = "{'HoursComputer': '1 - 4 hours', 'FormalEducation': 'Master’s degree (MA, MS, M.Eng., MBA, etc.)'}"
$ws.Range("D52").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 20000
$ws.Range("G52").Value = 0.6244888240034925
$ws.Range("H52").Value = 0.01040814706672488

$ws.Range("A53").Value = "2025-06-25 21:36:23"
$ws.Range("B53").Value = "Apriori"
$ws.Range("C53").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D53").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E53").Value = 25
$ws.Range("F53").Value = 20000
$ws.Range("G53").Value = 0.9656246330050635
$ws.Range("H53").Value = 0.01609374388341773

$ws.Range("A54").Value = "2025-06-25 21:36:25"
$ws.Range("B54").Value = "Apriori"
$ws.Range("C54").Value = "{'HoursComputer': '5 - 8 hours'}"
$ws.Range("D54").Value = "{'Gender': 'Male'}"
$ws.Range("E54").Value = 33
$ws.Range("F54").Value = 20000
$ws.Range("G54").Value = 1.457824452998466
$ws.Range("H54").Value = 0.0242970742166411

$ws.Range("A55").Value = "2025-06-25 21:36:27"
$ws.Range("B55").Value = "Apriori"
$ws.Range("C55").Value = "{'HoursComputer': '1 - 4 hours', 'FormalEducation': 'Master’s degree (MA, MS, M.Eng., MBA, etc.)'}"
$ws.Range("D55").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E55").Value = 43
$ws.Range("F55").Value = 20000
$ws.Range("G55").Value = 1.569658543005062
$ws.Range("H55").Value = 0.02616097571675103

$ws.Range("A56").Value = "2025-06-25 21:39:58"
$ws.Range("B56").Value = "Apriori"
$ws.Range("C56").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D56").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E56").Value = 21
$ws.Range("F56").Value = 20000
$ws.Range("G56").Value = 2.618635594000807
$ws.Range("H56").Value = 0.04364392656668012

$ws.Range("A57").Value = "2025-06-25 21:40:01"
$ws.Range("B57").Value = "Apriori"
$ws.Range("C57").Value = "{'HoursComputer': '5 - 8 hours'}"
$ws.Range("D57").Value = "{'Gender': 'Male'}"
$ws.Range("E57").Value = 29
$ws.Range("F57").Value = 20000
$ws.Range("G57").Value = 2.74235752900131
$ws.Range("H57").Value = 0.04570595881668851

$ws.Range("A58").Value = "2025-06-25 21:40:04"
$ws.Range("B58").Value = "Apriori"
$ws.Range("C58").Value = "{'HoursComputer': '1 - 4 hours', 'FormalEducation': 'Master’s degree (MA, MS, M.Eng., MBA, etc.)'}"
$ws.Range("D58").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E58").Value = 21
$ws.Range("F58").Value = 20000
$ws.Range("G58").Value = 1.559883781999815
$ws.Range("H58").Value = 0.02599806303333025

$ws.Range("A59").Value = "2025-06-25 21:54:16"
$ws.Range("B59").Value = "FP"
$ws.Range("C59").Value = "{'Exercise': '3 - 4 times per week'}"
$ws.Range("D59").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E59").Value = 21
$ws.Range("F59").Value = 20000
$ws.Range("G59").Value = 0.8526472130033653
$ws.Range("H59").Value = 0.01421078688338942

$ws.Range("A60").Value = "2025-06-25 21:54:18"
$ws.Range("B60").Value = "FP"
$ws.Range("C60").Value = "{'HoursComputer': '5 - 8 hours'}"
$ws.Range("D60").Value = "{'Gender': 'Male'}"
$ws.Range("E60").Value = 29
$ws.Range("F60").Value = 20000
$ws.Range("G60").Value = 1.825525783002377
$ws.Range("H60").Value = 0.03042542971670627

$ws.Range("A61").Value = "2025-06-25 21:54:20"
$ws.Range("B61").Value = "FP"
$ws.Range("C61").Value = "{'HoursComputer': '1 - 4 hours', 'FormalEducation': 'Master’s degree (MA, MS, M.Eng., MBA, etc.)'}"
$ws.Range("D61").Value = "{'RaceEthnicity': 'White or of European descent'}"
$ws.Range("E61").Value = "'21"
$ws.Range("F61").Value = "'20000"
$ws.Range("E61").ClearFormats()
$ws.Range("F61").ClearFormats()
$ws.Range("G61").Value = 1.416913036999176
$ws.Range("H61").Value = 0.0236152172833196

